$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") rows 2-35, replacing the old Strike# values
$newValues = @{
    2  = 3
    3  = 3
    4  = 3
    5  = 4
    6  = 4
    7  = 8
    8  = 6
    9  = 6
    10 = 7
    11 = 5
    12 = 1
    13 = 6
    14 = 6
    15 = 3
    16 = 2
    17 = 5
    18 = 3
    19 = 4
    20 = 2
    21 = 2
    22 = 3
    23 = 0
    24 = 3
    25 = 2
    26 = 1
    27 = 3
    28 = 3
    29 = 7
    30 = 2
    31 = 3
    32 = 0
    33 = 0
    34 = 1
    35 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
